$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.576.24'
$ws.Range("E2").Value = '  -0.06%  '

$ws.Range("D3").Value = '3.764.11'
$ws.Range("E3").Value = '  +0.57%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '616.25'
$ws.Range("E5").Value = '  +0.75%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '179.66'
$ws.Range("E6").Value = '  +1.06%  '

$ws.Range("D7").Value = '3.762.05'
$ws.Range("E7").Value = '  +0.63%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.527'
$ws.Range("E9").Value = '  -0.95%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.167'
$ws.Range("E10").Value = '  +1.20%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.55'
$ws.Range("E11").Value = '  +3.73%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.483'
$ws.Range("E12").Value = '  -1.72%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '40.17'
$ws.Range("E13").Value = '  -1.18%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000255'
$ws.Range("E14").Value = '  +1.12%  '

$ws.Range("D15").Value = '4.370.54'
$ws.Range("E15").Value = '  +0.04%  '

$ws.Range("D16").Value = '3.749.84'
$ws.Range("E16").Value = '  -0.09%  '

$ws.Range("D17").Value = '69.652.24'
$ws.Range("E17").Value = '  -0.02%  '

$ws.Range("E18").Value = '  -2.36%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.45'
$ws.Range("E19").Value = '  -1.21%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.44'
$ws.Range("E20").Value = '  -1.13%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '502.14'
$ws.Range("E21").Value = '  -2.00%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.30'
$ws.Range("E22").Value = '  -2.01%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.724'
$ws.Range("E23").Value = '  +0.22%  '

$ws.Range("E24").Value = '  +3.45%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.83'
$ws.Range("E25").Value = '  -2.11%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.97'
$ws.Range("E26").Value = '  -1.96%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.04'
$ws.Range("E27").Value = '  -0.02%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000136'
$ws.Range("E28").Value = '  +6.84%  '

$ws.Range("E29").Value = '  +0.16%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.52'
$ws.Range("E30").Value = '  +1.71%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.08'
$ws.Range("E31").Value = '  +3.70%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.92'
$ws.Range("E32").Value = '  +2.97%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '30.57'
$ws.Range("E33").Value = '  -2.49%  '

$ws.Range("E34").Value = '  -0.87%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.997'
$ws.Range("E35").Value = '  -0.14%  '

$ws.Range("B36").Value = 'Mantle'
$ws.Range("C36").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.05'
$ws.Range("E36").Value = '  +0.95%  '

$ws.Range("B37").Value = 'Filecoin'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.15'
$ws.Range("E37").Value = '  -0.58%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.351'
$ws.Range("E38").Value = '  +4.42%  '

$ws.Range("E39").Value = '  +4.43%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '467.78'
$ws.Range("E40").Value = '  +12.30%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.09'
$ws.Range("E41").Value = '  +13.70%  '

$ws.Range("E42").Value = '  -4.10%  '

$ws.Range("E43").Value = '  -2.81%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '45.00'
$ws.Range("E44").Value = '  +0.71%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.59'
$ws.Range("E45").Value = '  -1.97%  '

$ws.Range("D46").Value = '2.956.02'
$ws.Range("E46").Value = '  -3.47%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0361'
$ws.Range("E47").Value = '  -0.26%  '

$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.30'
$ws.Range("E48").Value = '  -1.27%  '

$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '139.17'
$ws.Range("E49").Value = '  +2.73%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.48'
$ws.Range("E51").Value = '  -0.99%  '
